# Update the "Price" (column D) values for the cryptocurrency rows whose
# quotes changed, per the "Updated symbol list" GitHub Actions commit.
#
# The Price column stores numeric-looking values as TEXT (inlineStr in the
# original workbook), e.g. "271.73". A plain $ws.Range(...).Value = "272.48"
# assignment would let Excel auto-detect the numeric-looking string and
# convert the cell to a real number (losing the text type and introducing
# floating point noise). Prefixing the value with a leading apostrophe is
# the standard Excel way of saying "treat this as text even though it looks
# like a number", so the cell keeps its original text semantics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = [ordered]@{
    2  = "272.48"
    3  = "23.13"
    4  = "6.355"
    5  = "0.06335"
    6  = "3.662"
    7  = "6.768"
    8  = "1.402"
    9  = "0.8364"
    10 = "0.1628"
    11 = "0.08392"
    12 = "0.03414"
    13 = "0.03158"
    14 = "0.09307"
    15 = "3.897"
    16 = "0.001723"
    17 = "0.04869"
    18 = "0.006204"
    19 = "0.005505"
    22 = "3.744"
    23 = "2.302"
    24 = "0.01388"
    25 = "0.3338"
    26 = "0.1217"
    27 = "0.0002674"
    40 = "0.04695"
    41 = "0.006899"
    42 = "0.1183"
    43 = "0.003318"
    44 = "0.01248"
    45 = "0.00006266"
    47 = "0.6978"
    48 = "0.1232"
    49 = "0.00002093"
}

foreach ($row in $priceUpdates.Keys) {
    $newPrice = $priceUpdates[$row]
    $ws.Range("D$row").Value = "'" + $newPrice
}
